$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 0.53
$ws.Range("C5").Value = 0.8100000000000001
$ws.Range("D5").Value = 0.68
$ws.Range("E5").Value = 0.93
$ws.Range("F5").Value = 0.54
